# Apply the "bi_generic_import / bi_import_chart_of_accounts" partner.xlsx edit:
#  - add two new columns (U: x_color, V: x_total_amount) to the Sample Partner sheet
#  - populate the two data rows with "red"/1500 and "blue"/2500
#  - widen the new column V to fit its header
#  - nudge the saved view state (tab ratio / scroll position / selection) to match
#    the resaved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new header cells -------------------------------------------------
$ws.Range("U1").Value = "x_color"
$ws.Range("V1").Value = "x_total_amount"

# --- new data cells -----------------------------------------------------
$ws.Range("U2").Value = "red"
$ws.Range("V2").Value = 1500

$ws.Range("U3").Value = "blue"
$ws.Range("V3").Value = 2500

# --- column widths for the newly introduced columns ---------------------
$ws.Columns.Item(21).ColumnWidth = 8.51
$ws.Columns.Item(22).ColumnWidth = 15.56

# --- view state (best effort - host/UI state) ---------------------------
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.5
$win.ScrollRow = 1
$win.ScrollColumn = 14

$ws.Range("U4").Select()
